# 06.07.19 Today Sales Details
# Apply an AutoFilter on column A ("DealerName") of the data range,
# keeping only rows whose value is "Tulip-2" visible (all other data
# rows become hidden), and update the sheet's current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The filtered data range (header row 3, data rows 4:58).
$filterRange = $ws.Range("A3:G58")

# Apply a standard "value" AutoFilter on the first column (DealerName),
# showing only rows where the value equals "Tulip-2". Using the
# xlFilterValues operator (7) with an array of accepted values produces
# the normal <filters><filter .../></filters> autoFilter markup and
# automatically hides the non-matching data rows.
$filterRange.AutoFilter(1, @("Tulip-2"), 7)

# Update the active selection/view as recorded after filtering.
$ws.Range("F71").Select()
